# Updated code for address
# Adds a new "alternate work location address" test scenario to the
# Employee_Details sheet: four new header columns (BJ:BM) and a new
# data row (row 7) beneath the existing scenarios.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employee_Details")

# New header columns for the alternate work location address fields.
$ws.Range("BJ1").Value = "altWorkLocationAddressLine1"
$ws.Range("BK1").Value = "altWorkLocationAddressType"
$ws.Range("BL1").Value = "altWorkLocationCountry"
$ws.Range("BM1").Value = "altWorkLocationZipCode"

# New scenario row.
$ws.Range("A7").Value = "NEW_ADDRESS_ALTERNATIVE_WORK_LOCATION"
$ws.Range("AN7").Value = "test2019"
$ws.Range("AX7").Value = 10171
$ws.Range("BJ7").Value = "Crescent Loop Dr"
$ws.Range("BK7").Value = "Alternate work Location Address"
$ws.Range("BL7").Value = "United States"
$ws.Range("BM7").Value = 48382

# Employee_Details becomes the active sheet/tab, with the new row's
# address-type cell selected.
$ws.Activate() | Out-Null
$ws.Range("BK7").Select() | Out-Null
